$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B (TB), C (d2S), D (K), E (IP), G (sum) values for rows 2-49,
# regenerated after filtering save games out of the underlying per-game data.
$data = @(
    @{Row=2; B=3.272327238179451; C=9.98352242611593; D=18.71679738969934; E=13.86384647080068; G=45.8364935247954},
    @{Row=3; B=1.445647641019636; C=0.3048912486333797; D=3.223369029078222; E=0.5333859586016987; G=5.507293877332936},
    @{Row=4; B=3.272327238179451; C=1.626987699542094; D=0.1496068669990043; E=0.5333859586016987; G=5.582307763322248},
    @{Row=5; B=1.445647641019636; C=0.3048912486333797; D=0.7210945179870265; E=0.5333859586016987; G=3.005019366241741},
    @{Row=6; B=0.04172184405617529; C=1.626987699542094; D=0.7210945179870265; E=0.5333859586016987; G=2.923190020186994},
    @{Row=7; B=1.445647641019636; C=1.626987699542094; D=0.1496068669990043; E=0.5333859586016987; G=3.755628166162433},
    @{Row=8; B=1.445647641019636; C=1.626987699542094; D=0.7210945179870265; E=0.5333859586016987; G=4.327115817150455},
    @{Row=9; B=0.01253208636536152; C=0.002658071450198252; D=0.7210945179870265; E=0.5333859586016987; G=1.269670634404285},
    @{Row=10; B=1.445647641019636; C=1.626987699542094; D=3.223369029078222; E=0.5333859586016987; G=6.82939032824165},
    @{Row=11; B=1.445647641019636; C=1.626987699542094; D=0.7210945179870265; E=0.5333859586016987; G=4.327115817150455},
    @{Row=12; B=3.272327238179451; C=1.626987699542094; D=0.7210945179870265; E=0.5333859586016987; G=6.15379541431027},
    @{Row=13; B=1.445647641019636; C=1.626987699542094; D=0.7210945179870265; E=0.5333859586016987; G=4.327115817150455},
    @{Row=14; B=0.2881169905109251; C=0.3048912486333797; D=0.1496068669990043; E=0.5333859586016987; G=1.276001064745008},
    @{Row=15; B=0.6545652718822623; C=0.3048912486333797; D=18.71679738969934; E=0.5333859586016987; G=20.20963986881668},
    @{Row=16; B=3.272327238179451; C=1.626987699542094; D=0.7210945179870265; E=0.5333859586016987; G=6.15379541431027},
    @{Row=17; B=3.272327238179451; C=1.626987699542094; D=0.1496068669990043; E=0.5333859586016987; G=5.582307763322248},
    @{Row=18; B=1.445647641019636; C=1.626987699542094; D=3.223369029078222; E=0.5333859586016987; G=6.82939032824165},
    @{Row=19; B=3.272327238179451; C=1.626987699542094; D=0.1496068669990043; E=0.5333859586016987; G=5.582307763322248},
    @{Row=20; B=1.445647641019636; C=1.626987699542094; D=0.1496068669990043; E=0.5333859586016987; G=3.755628166162433},
    @{Row=21; B=0.6545652718822623; C=1.626987699542094; D=3.223369029078222; E=0.5333859586016987; G=6.038307959104277},
    @{Row=22; B=3.272327238179451; C=1.626987699542094; D=3.223369029078222; E=0.5333859586016987; G=8.656069925401464},
    @{Row=23; B=0.6545652718822623; C=1.626987699542094; D=3.223369029078222; E=0.5333859586016987; G=6.038307959104277},
    @{Row=24; B=1.445647641019636; C=1.626987699542094; D=0.1496068669990043; E=0.5333859586016987; G=3.755628166162433},
    @{Row=25; B=0.1169995834814548; C=0.04103571897497393; D=18.71679738969934; E=0.5333859586016987; G=19.40821865075747},
    @{Row=26; B=3.272327238179451; C=1.626987699542094; D=0.7210945179870265; E=0.5333859586016987; G=6.15379541431027},
    @{Row=27; B=0.2881169905109251; C=0.3048912486333797; D=0.1496068669990043; E=0.5333859586016987; G=1.276001064745008},
    @{Row=28; B=0.1169995834814548; C=0.3048912486333797; D=0.7210945179870265; E=0.5333859586016987; G=1.67637130870356},
    @{Row=29; B=1.445647641019636; C=1.626987699542094; D=0.1496068669990043; E=0.5333859586016987; G=3.755628166162433},
    @{Row=30; B=3.272327238179451; C=1.626987699542094; D=0.1496068669990043; E=0.5333859586016987; G=5.582307763322248},
    @{Row=31; B=0.01253208636536152; C=0.04103571897497393; D=0.7210945179870265; E=0.5333859586016987; G=1.308048281929061},
    @{Row=32; B=1.445647641019636; C=1.626987699542094; D=0.7210945179870265; E=0.5333859586016987; G=4.327115817150455},
    @{Row=33; B=0.6545652718822623; C=9.98352242611593; D=3.223369029078222; E=13.86384647080068; G=27.7253031978771},
    @{Row=34; B=3.272327238179451; C=1.626987699542094; D=3.223369029078222; E=0.5333859586016987; G=8.656069925401464},
    @{Row=35; B=3.272327238179451; C=1.626987699542094; D=18.71679738969934; E=0.5333859586016987; G=24.14949828602258},
    @{Row=36; B=3.272327238179451; C=109.9114832445916; D=18.71679738969934; E=13.86384647080068; G=145.7644543432711},
    @{Row=37; B=0.6545652718822623; C=1.626987699542094; D=18.71679738969934; E=0.5333859586016987; G=21.53173631972539},
    @{Row=38; B=0.6545652718822623; C=1.626987699542094; D=3.223369029078222; E=0.5333859586016987; G=6.038307959104277},
    @{Row=39; B=0.6545652718822623; C=1.626987699542094; D=0.7210945179870265; E=0.5333859586016987; G=3.536033448013082},
    @{Row=40; B=3.272327238179451; C=1.626987699542094; D=0.1496068669990043; E=0.5333859586016987; G=5.582307763322248},
    @{Row=41; B=3.272327238179451; C=0.3048912486333797; D=3.223369029078222; E=0.5333859586016987; G=7.333973474492751},
    @{Row=42; B=0.1169995834814548; C=0.3048912486333797; D=0.1496068669990043; E=0.5333859586016987; G=1.104883657715537},
    @{Row=43; B=0.6545652718822623; C=0.3048912486333797; D=0.1496068669990043; E=0.5333859586016987; G=1.642449346116345},
    @{Row=44; B=0.6545652718822623; C=1.626987699542094; D=3.223369029078222; E=0.5333859586016987; G=6.038307959104277},
    @{Row=45; B=3.272327238179451; C=1.626987699542094; D=0.7210945179870265; E=0.5333859586016987; G=6.15379541431027},
    @{Row=46; B=0.6545652718822623; C=1.626987699542094; D=3.223369029078222; E=0.5333859586016987; G=6.038307959104277},
    @{Row=47; B=0.2881169905109251; C=1.626987699542094; D=3.223369029078222; E=0.5333859586016987; G=5.671859677732939},
    @{Row=48; B=3.272327238179451; C=1.626987699542094; D=0.7210945179870265; E=0.5333859586016987; G=6.15379541431027},
    @{Row=49; B=3.272327238179451; C=1.626987699542094; D=3.223369029078222; E=0.5333859586016987; G=6.15379541431027}
)

foreach ($item in $data) {
    $ws.Range("B" + $item.Row).Value = $item.B
    $ws.Range("C" + $item.Row).Value = $item.C
    $ws.Range("D" + $item.Row).Value = $item.D
    $ws.Range("E" + $item.Row).Value = $item.E
    $ws.Range("G" + $item.Row).Value = $item.G
}
